$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129; this shifts existing rows 129-194 down
# to 130-195 (carrying their values/formatting with them, including the
# last row's data landing on the new row 195).
$ws.Rows(129).Insert()

# Populate the newly inserted row 129 with this week's new record. All the
# "descriptive" columns (market, region, product taxonomy, unit, etc.) are
# constant for every row in this sheet, so mirror them here; only the
# date/volume/price/origin columns carry genuinely new data.
$ws.Range("A129").Value = 1
$ws.Range("B129").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C129").Value = "Arica y Parinacota"
$ws.Range("D129").Value = 44606
$ws.Range("E129").Value = 15
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100108
$ws.Range("H129").Value = "Tropicales y subtropicales"
$ws.Range("I129").Value = 100108006
$ws.Range("J129").Value = "Plátano"
$ws.Range("K129").Value = "Sin especificar"
$ws.Range("L129").Value = "Pintón"
$ws.Range("M129").Value = 120
$ws.Range("N129").Value = 17000
$ws.Range("O129").Value = 18000
$ws.Range("P129").Value = 17500
$ws.Range("Q129").Value = "$/caja 20 kilos"
$ws.Range("R129").Value = "Ecuador"
$ws.Range("S129").Value = 875
$ws.Range("T129").Value = 20
